$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 28, pushing existing rows 28-61 down to 29-62.
$ws.Rows.Item(28).Insert()

# Fill the new row 28 with the same data as the (now shifted) row 29,
# except for the date (one day later) and the updated volume.
$ws.Range("A28:R28").Value2 = $ws.Range("A29:R29").Value2

$ws.Range("D28").Value2 = 44587
$ws.Range("J28").Value2 = 400
